$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-10-25 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-26 Saturday", 2) | Out-Null

# Update each arithmetic expression cell (MatchWholeWord avoids substring collisions,
# e.g. "9+26=" vs "29+26=", and "0+24=" vs "50+24=")
$d.Content.Find.Execute("59-2=", $true, $true, $false, $false, $false, $true, 1, $false, "66-9=", 2) | Out-Null
$d.Content.Find.Execute("72+18=", $true, $true, $false, $false, $false, $true, 1, $false, "31-3=", 2) | Out-Null
$d.Content.Find.Execute("57-49=", $true, $true, $false, $false, $false, $true, 1, $false, "51+7=", 2) | Out-Null
$d.Content.Find.Execute("73-25=", $true, $true, $false, $false, $false, $true, 1, $false, "15+33=", 2) | Out-Null
$d.Content.Find.Execute("47-30=", $true, $true, $false, $false, $false, $true, 1, $false, "44+0=", 2) | Out-Null
$d.Content.Find.Execute("64-21=", $true, $true, $false, $false, $false, $true, 1, $false, "55+7=", 2) | Out-Null
$d.Content.Find.Execute("62-53=", $true, $true, $false, $false, $false, $true, 1, $false, "15+70=", 2) | Out-Null
$d.Content.Find.Execute("73+9=", $true, $true, $false, $false, $false, $true, 1, $false, "81-64=", 2) | Out-Null
$d.Content.Find.Execute("98-15=", $true, $true, $false, $false, $false, $true, 1, $false, "45+23=", 2) | Out-Null
$d.Content.Find.Execute("11+25=", $true, $true, $false, $false, $false, $true, 1, $false, "71+21=", 2) | Out-Null
$d.Content.Find.Execute("73+18=", $true, $true, $false, $false, $false, $true, 1, $false, "37+18=", 2) | Out-Null
$d.Content.Find.Execute("39-38=", $true, $true, $false, $false, $false, $true, 1, $false, "91-46=", 2) | Out-Null
$d.Content.Find.Execute("59-3=", $true, $true, $false, $false, $false, $true, 1, $false, "34+3=", 2) | Out-Null
$d.Content.Find.Execute("74-56=", $true, $true, $false, $false, $false, $true, 1, $false, "0+24=", 2) | Out-Null
$d.Content.Find.Execute("9+24=", $true, $true, $false, $false, $false, $true, 1, $false, "18+8=", 2) | Out-Null
$d.Content.Find.Execute("71-10=", $true, $true, $false, $false, $false, $true, 1, $false, "41+7=", 2) | Out-Null
$d.Content.Find.Execute("21+11=", $true, $true, $false, $false, $false, $true, 1, $false, "50+31=", 2) | Out-Null
$d.Content.Find.Execute("31+16=", $true, $true, $false, $false, $false, $true, 1, $false, "81-54=", 2) | Out-Null
$d.Content.Find.Execute("51-31=", $true, $true, $false, $false, $false, $true, 1, $false, "54-19=", 2) | Out-Null
$d.Content.Find.Execute("8+78=", $true, $true, $false, $false, $false, $true, 1, $false, "78-63=", 2) | Out-Null
$d.Content.Find.Execute("7+61=", $true, $true, $false, $false, $false, $true, 1, $false, "86-28=", 2) | Out-Null
$d.Content.Find.Execute("52-4=", $true, $true, $false, $false, $false, $true, 1, $false, "78-54=", 2) | Out-Null
$d.Content.Find.Execute("33-18=", $true, $true, $false, $false, $false, $true, 1, $false, "12+55=", 2) | Out-Null
$d.Content.Find.Execute("39+20=", $true, $true, $false, $false, $false, $true, 1, $false, "3+54=", 2) | Out-Null
$d.Content.Find.Execute("74-69=", $true, $true, $false, $false, $false, $true, 1, $false, "20-7=", 2) | Out-Null
$d.Content.Find.Execute("54+13=", $true, $true, $false, $false, $false, $true, 1, $false, "57+1=", 2) | Out-Null
$d.Content.Find.Execute("21+68=", $true, $true, $false, $false, $false, $true, 1, $false, "79-37=", 2) | Out-Null
$d.Content.Find.Execute("35+44=", $true, $true, $false, $false, $false, $true, 1, $false, "49-2=", 2) | Out-Null
$d.Content.Find.Execute("60-13=", $true, $true, $false, $false, $false, $true, 1, $false, "98-64=", 2) | Out-Null
$d.Content.Find.Execute("50-42=", $true, $true, $false, $false, $false, $true, 1, $false, "22+37=", 2) | Out-Null
$d.Content.Find.Execute("9+29=", $true, $true, $false, $false, $false, $true, 1, $false, "56+42=", 2) | Out-Null
$d.Content.Find.Execute("13+8=", $true, $true, $false, $false, $false, $true, 1, $false, "53-27=", 2) | Out-Null
$d.Content.Find.Execute("37-22=", $true, $true, $false, $false, $false, $true, 1, $false, "53+4=", 2) | Out-Null
$d.Content.Find.Execute("54+5=", $true, $true, $false, $false, $false, $true, 1, $false, "93-14=", 2) | Out-Null
$d.Content.Find.Execute("92-36=", $true, $true, $false, $false, $false, $true, 1, $false, "83-57=", 2) | Out-Null
$d.Content.Find.Execute("78+3=", $true, $true, $false, $false, $false, $true, 1, $false, "41+26=", 2) | Out-Null
$d.Content.Find.Execute("77-3=", $true, $true, $false, $false, $false, $true, 1, $false, "31-10=", 2) | Out-Null
$d.Content.Find.Execute("17+8=", $true, $true, $false, $false, $false, $true, 1, $false, "85-81=", 2) | Out-Null
$d.Content.Find.Execute("40-23=", $true, $true, $false, $false, $false, $true, 1, $false, "19-2=", 2) | Out-Null
$d.Content.Find.Execute("89-63=", $true, $true, $false, $false, $false, $true, 1, $false, "1+95=", 2) | Out-Null
$d.Content.Find.Execute("9+26=", $true, $true, $false, $false, $false, $true, 1, $false, "82-47=", 2) | Out-Null
$d.Content.Find.Execute("68+18=", $true, $true, $false, $false, $false, $true, 1, $false, "1+47=", 2) | Out-Null
$d.Content.Find.Execute("65+11=", $true, $true, $false, $false, $false, $true, 1, $false, "97-48=", 2) | Out-Null
$d.Content.Find.Execute("31+41=", $true, $true, $false, $false, $false, $true, 1, $false, "90+4=", 2) | Out-Null
$d.Content.Find.Execute("27+69=", $true, $true, $false, $false, $false, $true, 1, $false, "50+24=", 2) | Out-Null
$d.Content.Find.Execute("26+5=", $true, $true, $false, $false, $false, $true, 1, $false, "43-38=", 2) | Out-Null
$d.Content.Find.Execute("18-11=", $true, $true, $false, $false, $false, $true, 1, $false, "37+22=", 2) | Out-Null
$d.Content.Find.Execute("6+71=", $true, $true, $false, $false, $false, $true, 1, $false, "25+62=", 2) | Out-Null
$d.Content.Find.Execute("41-22=", $true, $true, $false, $false, $false, $true, 1, $false, "39-0=", 2) | Out-Null
$d.Content.Find.Execute("86+8=", $true, $true, $false, $false, $false, $true, 1, $false, "62-55=", 2) | Out-Null
$d.Content.Find.Execute("68-0=", $true, $true, $false, $false, $false, $true, 1, $false, "0+62=", 2) | Out-Null
$d.Content.Find.Execute("50-3=", $true, $true, $false, $false, $false, $true, 1, $false, "1+91=", 2) | Out-Null
$d.Content.Find.Execute("6+1=", $true, $true, $false, $false, $false, $true, 1, $false, "69-10=", 2) | Out-Null
$d.Content.Find.Execute("26+16=", $true, $true, $false, $false, $false, $true, 1, $false, "45+49=", 2) | Out-Null
$d.Content.Find.Execute("32+42=", $true, $true, $false, $false, $false, $true, 1, $false, "27+32=", 2) | Out-Null
$d.Content.Find.Execute("11-3=", $true, $true, $false, $false, $false, $true, 1, $false, "9+87=", 2) | Out-Null
$d.Content.Find.Execute("30+33=", $true, $true, $false, $false, $false, $true, 1, $false, "87-5=", 2) | Out-Null
$d.Content.Find.Execute("8+64=", $true, $true, $false, $false, $false, $true, 1, $false, "50-30=", 2) | Out-Null
$d.Content.Find.Execute("11+5=", $true, $true, $false, $false, $false, $true, 1, $false, "54-37=", 2) | Out-Null
$d.Content.Find.Execute("14+27=", $true, $true, $false, $false, $false, $true, 1, $false, "5-3=", 2) | Out-Null
$d.Content.Find.Execute("80+7=", $true, $true, $false, $false, $false, $true, 1, $false, "70-5=", 2) | Out-Null
$d.Content.Find.Execute("39-37=", $true, $true, $false, $false, $false, $true, 1, $false, "20-3=", 2) | Out-Null
$d.Content.Find.Execute("14+61=", $true, $true, $false, $false, $false, $true, 1, $false, "16-3=", 2) | Out-Null
$d.Content.Find.Execute("30-18=", $true, $true, $false, $false, $false, $true, 1, $false, "70-45=", 2) | Out-Null
$d.Content.Find.Execute("45-26=", $true, $true, $false, $false, $false, $true, 1, $false, "12+79=", 2) | Out-Null
$d.Content.Find.Execute("39-28=", $true, $true, $false, $false, $false, $true, 1, $false, "8+29=", 2) | Out-Null
$d.Content.Find.Execute("80-3=", $true, $true, $false, $false, $false, $true, 1, $false, "93-16=", 2) | Out-Null
$d.Content.Find.Execute("41-31=", $true, $true, $false, $false, $false, $true, 1, $false, "85-7=", 2) | Out-Null
$d.Content.Find.Execute("3+20=", $true, $true, $false, $false, $false, $true, 1, $false, "97-85=", 2) | Out-Null
$d.Content.Find.Execute("32+39=", $true, $true, $false, $false, $false, $true, 1, $false, "71-70=", 2) | Out-Null
$d.Content.Find.Execute("28+4=", $true, $true, $false, $false, $false, $true, 1, $false, "26+11=", 2) | Out-Null
$d.Content.Find.Execute("47+24=", $true, $true, $false, $false, $false, $true, 1, $false, "94-57=", 2) | Out-Null
$d.Content.Find.Execute("75-61=", $true, $true, $false, $false, $false, $true, 1, $false, "57-1=", 2) | Out-Null
$d.Content.Find.Execute("40+49=", $true, $true, $false, $false, $false, $true, 1, $false, "23-19=", 2) | Out-Null
$d.Content.Find.Execute("91-26=", $true, $true, $false, $false, $false, $true, 1, $false, "49-35=", 2) | Out-Null
$d.Content.Find.Execute("2+90=", $true, $true, $false, $false, $false, $true, 1, $false, "54+26=", 2) | Out-Null
$d.Content.Find.Execute("98-80=", $true, $true, $false, $false, $false, $true, 1, $false, "9+59=", 2) | Out-Null
$d.Content.Find.Execute("11+44=", $true, $true, $false, $false, $false, $true, 1, $false, "41-21=", 2) | Out-Null
$d.Content.Find.Execute("92-24=", $true, $true, $false, $false, $false, $true, 1, $false, "28+13=", 2) | Out-Null
$d.Content.Find.Execute("65-40=", $true, $true, $false, $false, $false, $true, 1, $false, "81-21=", 2) | Out-Null
$d.Content.Find.Execute("1+94=", $true, $true, $false, $false, $false, $true, 1, $false, "15+22=", 2) | Out-Null
$d.Content.Find.Execute("87+7=", $true, $true, $false, $false, $false, $true, 1, $false, "49+44=", 2) | Out-Null
$d.Content.Find.Execute("38-30=", $true, $true, $false, $false, $false, $true, 1, $false, "36+27=", 2) | Out-Null
$d.Content.Find.Execute("97-86=", $true, $true, $false, $false, $false, $true, 1, $false, "27+43=", 2) | Out-Null
$d.Content.Find.Execute("92-37=", $true, $true, $false, $false, $false, $true, 1, $false, "51-14=", 2) | Out-Null
$d.Content.Find.Execute("21+22=", $true, $true, $false, $false, $false, $true, 1, $false, "27+36=", 2) | Out-Null
$d.Content.Find.Execute("28+15=", $true, $true, $false, $false, $false, $true, 1, $false, "32+33=", 2) | Out-Null
$d.Content.Find.Execute("99-51=", $true, $true, $false, $false, $false, $true, 1, $false, "40+38=", 2) | Out-Null
$d.Content.Find.Execute("40-37=", $true, $true, $false, $false, $false, $true, 1, $false, "68-66=", 2) | Out-Null
$d.Content.Find.Execute("18-15=", $true, $true, $false, $false, $false, $true, 1, $false, "53-3=", 2) | Out-Null
$d.Content.Find.Execute("54+8=", $true, $true, $false, $false, $false, $true, 1, $false, "4+24=", 2) | Out-Null
$d.Content.Find.Execute("28-20=", $true, $true, $false, $false, $false, $true, 1, $false, "93+3=", 2) | Out-Null
$d.Content.Find.Execute("68-59=", $true, $true, $false, $false, $false, $true, 1, $false, "49+5=", 2) | Out-Null
$d.Content.Find.Execute("36+6=", $true, $true, $false, $false, $false, $true, 1, $false, "38+39=", 2) | Out-Null
$d.Content.Find.Execute("65+3=", $true, $true, $false, $false, $false, $true, 1, $false, "42-3=", 2) | Out-Null
$d.Content.Find.Execute("68+10=", $true, $true, $false, $false, $false, $true, 1, $false, "32-28=", 2) | Out-Null
$d.Content.Find.Execute("17+56=", $true, $true, $false, $false, $false, $true, 1, $false, "89-69=", 2) | Out-Null
$d.Content.Find.Execute("14+21=", $true, $true, $false, $false, $false, $true, 1, $false, "18+31=", 2) | Out-Null
$d.Content.Find.Execute("29+26=", $true, $true, $false, $false, $false, $true, 1, $false, "89-10=", 2) | Out-Null
$d.Content.Find.Execute("56+11=", $true, $true, $false, $false, $false, $true, 1, $false, "87-49=", 2) | Out-Null
